# Update the five "two-digit division" practice rows with new problems.
# Each data row lives at 1-based table row 1, 5, 9, 13, 17; the other rows
# are blank answer rows. We overwrite each cell's Range.Text directly so
# the existing run/paragraph formatting (font, size, justification) is
# preserved exactly as in the original document.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "66÷8="
$t.Cell(1, 2).Range.Text = "44÷6="
$t.Cell(1, 3).Range.Text = "94÷9="
$t.Cell(1, 4).Range.Text = "92÷3="
$t.Cell(1, 5).Range.Text = "28÷2="

$t.Cell(5, 1).Range.Text = "43÷9="
$t.Cell(5, 2).Range.Text = "76÷3="
$t.Cell(5, 3).Range.Text = "48÷6="
$t.Cell(5, 4).Range.Text = "54÷9="
$t.Cell(5, 5).Range.Text = "44÷3="

$t.Cell(9, 1).Range.Text = "10÷7="
$t.Cell(9, 2).Range.Text = "37÷9="
$t.Cell(9, 3).Range.Text = "21÷5="
$t.Cell(9, 4).Range.Text = "10÷5="
$t.Cell(9, 5).Range.Text = "54÷2="

$t.Cell(13, 1).Range.Text = "41÷9="
$t.Cell(13, 2).Range.Text = "88÷3="
$t.Cell(13, 3).Range.Text = "13÷3="
$t.Cell(13, 4).Range.Text = "82÷5="
$t.Cell(13, 5).Range.Text = "51÷6="

$t.Cell(17, 1).Range.Text = "96÷2="
$t.Cell(17, 2).Range.Text = "86÷6="
$t.Cell(17, 3).Range.Text = "46÷2="
$t.Cell(17, 4).Range.Text = "49÷5="
$t.Cell(17, 5).Range.Text = "12÷8="
